# Scheduled-runner market data refresh: updates the per-leve price/profit
# columns (H:N) on several job sheets. Cells that had no prior value get a
# new one, and a few cells that no longer apply are cleared outright.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 1039
$ws.Range("J34").Value = 5
$ws.Range("L34").Value = 5
$ws.Range("N34").Value = -411
$ws.Range("H36").Value = 1039
$ws.Range("J36").Value = 5
$ws.Range("L36").Value = 5
$ws.Range("N36").Value = -1435
$ws.Range("H43").Value = 302
$ws.Range("J43").Value = 302
$ws.Range("L43").Value = 302
$ws.Range("N43").Value = -440
$ws.Range("H58").Value = 478
$ws.Range("I58").Value = 130
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 390
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -240
$ws.Range("N58").Value = -3300
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()
$ws.Range("H113").Value = 2802.5
$ws.Range("I113").Value = 2524.2856
$ws.Range("K113").Value = 2524.2856
$ws.Range("M113").Value = 729.7143999999998
$ws.Range("H129").Value = 20001600
$ws.Range("I129").Value = 50000400
$ws.Range("J129").Value = 2400
$ws.Range("K129").Value = 150001200
$ws.Range("L129").Value = 7200
$ws.Range("M129").Value = -149996200
$ws.Range("N129").Value = -17200
$ws.Range("H138").Value = 10003840
$ws.Range("J138").Value = 3982.6667
$ws.Range("L138").Value = 11948.0001
$ws.Range("N138").Value = -22228.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5724.8335
$ws.Range("I63").Value = 5869.8
$ws.Range("K63").Value = 5869.8
$ws.Range("M63").Value = -5183.8
$ws.Range("H66").Value = 5724.8335
$ws.Range("I66").Value = 5869.8
$ws.Range("K66").Value = 29349
$ws.Range("M66").Value = -25917
$ws.Range("H140").Value = 52762
$ws.Range("I140").Value = 55555
$ws.Range("J140").Value = 49969
$ws.Range("K140").Value = 55555
$ws.Range("L140").Value = 49969
$ws.Range("M140").Value = -50375
$ws.Range("N140").Value = -60329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2350.25
$ws.Range("I105").Value = 2300.3333
$ws.Range("K105").Value = 2300.3333
$ws.Range("M105").Value = -553.3332999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8427.375
$ws.Range("I58").Value = 6570.1665
$ws.Range("K58").Value = 6570.1665
$ws.Range("M58").Value = -6367.1665
$ws.Range("H86").Value = 13071.429
$ws.Range("I86").Value = 12500
$ws.Range("J86").Value = 14500
$ws.Range("K86").Value = 12500
$ws.Range("L86").Value = 14500
$ws.Range("M86").Value = -11377
$ws.Range("N86").Value = -16746
$ws.Range("H89").Value = 13071.429
$ws.Range("I89").Value = 12500
$ws.Range("J89").Value = 14500
$ws.Range("K89").Value = 62500
$ws.Range("L89").Value = 72500
$ws.Range("M89").Value = -56884
$ws.Range("N89").Value = -83732
$ws.Range("H105").Value = 2404.5
$ws.Range("I105").Value = 2085.4
$ws.Range("K105").Value = 2085.4
$ws.Range("M105").Value = -338.4000000000001
$ws.Range("H136").Value = 8427.375
$ws.Range("I136").Value = 6570.1665
$ws.Range("K136").Value = 19710.4995
$ws.Range("M136").Value = -17160.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 185
$ws.Range("I10").Value = 80
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 240
$ws.Range("L10").Value = 1500
$ws.Range("M10").Value = -101
$ws.Range("N10").Value = -1778
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H80").Value = 1495
$ws.Range("J80").Value = 1495
$ws.Range("L80").Value = 4485
$ws.Range("N80").Value = -6357
$ws.Range("H83").Value = 1495
$ws.Range("J83").Value = 1495
$ws.Range("L83").Value = 13455
$ws.Range("N83").Value = -22815
$ws.Range("H104").Value = 5992.3335
$ws.Range("J104").Value = 5992.3335
$ws.Range("L104").Value = 17977.0005
$ws.Range("N104").Value = -23219.0005
$ws.Range("H117").Value = 3210.0908
$ws.Range("J117").Value = 3031.1
$ws.Range("L117").Value = 9093.299999999999
$ws.Range("N117").Value = -15977.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2981.25
$ws.Range("I3").Value = 3000
$ws.Range("J3").Value = 2975
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 2975
$ws.Range("M3").Value = -2884
$ws.Range("N3").Value = -3207
$ws.Range("H43").Value = 16338.333
$ws.Range("I43").Value = 14508
$ws.Range("K43").Value = 14508
$ws.Range("M43").Value = -14357
$ws.Range("H70").Value = 5715.857
$ws.Range("I70").Value = 5668.6665
$ws.Range("J70").Value = 5999
$ws.Range("K70").Value = 5668.6665
$ws.Range("L70").Value = 5999
$ws.Range("M70").Value = -5398.6665
$ws.Range("N70").Value = -6539
$ws.Range("H73").Value = 5715.857
$ws.Range("I73").Value = 5668.6665
$ws.Range("J73").Value = 5999
$ws.Range("K73").Value = 5668.6665
$ws.Range("L73").Value = 5999
$ws.Range("M73").Value = -4732.6665
$ws.Range("N73").Value = -7871
$ws.Range("H80").Value = 4000
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 4000
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 102750
$ws.Range("J62").Value = 102750
$ws.Range("L62").Value = 102750
$ws.Range("N62").Value = -103998
$ws.Range("H65").Value = 102750
$ws.Range("J65").Value = 102750
$ws.Range("L65").Value = 513750
$ws.Range("N65").Value = -519990
$ws.Range("H75").Value = 25000
$ws.Range("I75").Value = 20000
$ws.Range("K75").Value = 20000
$ws.Range("M75").Value = -19064
$ws.Range("H78").Value = 25000
$ws.Range("I78").Value = 20000
$ws.Range("K78").Value = 60000
$ws.Range("M78").Value = -55320
